$d = $word.ActiveDocument

# --- 1. Remove the "Meta description" paragraph (paragraph 2) ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Meta description")) {
        $p.Range.Delete()
        break
    }
}

# --- 2. Split the last paragraph (the italic image-prompt paragraph) into two ---
#     a) a new bold paragraph: "Play Cleopatra Diamond Spins Free | Review of IGT's Slot Game"
#     b) the existing paragraph's text replaced with the meta-description sentence (keeping italics)
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$prevPara = $d.Paragraphs.Item($count - 1)

$insertionPoint = $prevPara.Range.Duplicate
$insertionPoint.Collapse(0)
$newHeadingText = "Play Cleopatra Diamond Spins Free | Review of IGT's Slot Game`r"
$insertionPoint.InsertAfter($newHeadingText)

$newPara = $d.Paragraphs.Item($count)
$newParaTextRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newParaTextRange.Font.Bold = $true

# Now update the text of the (shifted) last paragraph, preserving its italic run formatting
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$find = $lastPara.Range.Find
$find.ClearFormatting()
$find.Text = "Please create a feature image for ""Cleopatra Diamond Spins"" that meets the following criteria: - Cartoon style - Features a happy Maya warrior with glasses. The image should convey a sense of fun and excitement, while also incorporating the historical theme of the game. The Maya warrior should be portrayed in a way that suggests they are enjoying playing the game, with their glasses indicating that they are knowledgeable and experienced players. The color scheme should be vibrant and eye-catching, and there should be some reference to the Ancient Egyptian setting of the game, such as hieroglyphics or pyramids in the background."
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "Learn about IGT's Cleopatra Diamond Spins slot game in this review, and play for free. Features, symbols, gameplay, and sound and graphics are discussed."
$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null
